$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(35, 1).Value = "^_^^_^"
$ws.Cells.Item(35, 2).Value = "Verificar"
$ws.Cells.Item(35, 3).Value = "Guerra Atual"
$ws.Cells.Item(35, 4).Value = 8
$ws.Cells.Item(35, 5).Value = 12
$ws.Cells.Item(35, 6).Value = 16
$ws.Cells.Item(35, 7).Value = 12
$ws.Cells.Item(35, 8).Value = 12

$ws.Cells.Item(36, 1).Value = "ZackThunder"
$ws.Cells.Item(36, 2).Value = "Ok"
$ws.Cells.Item(36, 3).Value = "Guerra Atual"
$ws.Cells.Item(36, 4).Value = 16
$ws.Cells.Item(36, 5).Value = 16
$ws.Cells.Item(36, 6).Value = 8
$ws.Cells.Item(36, 7).Value = 16
$ws.Cells.Item(36, 8).Value = 16

$ws.Cells.Item(37, 1).Value = "polaris"
$ws.Cells.Item(37, 2).Value = "Ok"
$ws.Cells.Item(37, 3).Value = "Guerra Atual"
$ws.Cells.Item(37, 4).Value = 16
$ws.Cells.Item(37, 5).Value = 16
$ws.Cells.Item(37, 6).Value = 16
$ws.Cells.Item(37, 7).Value = 16
$ws.Cells.Item(37, 8).Value = 16

$ws.Cells.Item(38, 1).Value = "Gustavo Clash"
$ws.Cells.Item(38, 2).Value = "Ok"
$ws.Cells.Item(38, 3).Value = "Guerra Atual"
$ws.Cells.Item(38, 4).Value = 16
$ws.Cells.Item(38, 5).Value = 16
$ws.Cells.Item(38, 6).Value = 16
$ws.Cells.Item(38, 7).Value = 16
$ws.Cells.Item(38, 8).Value = 16

$ws.Cells.Item(39, 1).Value = "RaiNascimento"
$ws.Cells.Item(39, 2).Value = "Ok"
$ws.Cells.Item(39, 3).Value = "Guerra Atual"
$ws.Cells.Item(39, 4).Value = 16
$ws.Cells.Item(39, 5).Value = 0
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(39, 7).Value = 0
$ws.Cells.Item(39, 8).Value = 0

$ws.Cells.Item(40, 1).Value = "super"
$ws.Cells.Item(40, 2).Value = "Ok"
$ws.Cells.Item(40, 3).Value = "Guerra Atual"
$ws.Cells.Item(40, 4).Value = 16
$ws.Cells.Item(40, 5).Value = 16
$ws.Cells.Item(40, 6).Value = 16
$ws.Cells.Item(40, 7).Value = 16
$ws.Cells.Item(40, 8).Value = 16

$ws.Cells.Item(41, 1).Value = "RobaFrag"
$ws.Cells.Item(41, 2).Value = "Ok"
$ws.Cells.Item(41, 3).Value = "Guerra Atual"
$ws.Cells.Item(41, 4).Value = 16
$ws.Cells.Item(41, 5).Value = 16
$ws.Cells.Item(41, 6).Value = 16
$ws.Cells.Item(41, 7).Value = 16
$ws.Cells.Item(41, 8).Value = 16

$ws.Cells.Item(42, 1).Value = "Rodolfos"
$ws.Cells.Item(42, 2).Value = "Ok"
$ws.Cells.Item(42, 3).Value = "Guerra Atual"
$ws.Cells.Item(42, 4).Value = 16
$ws.Cells.Item(42, 5).Value = 16
$ws.Cells.Item(42, 6).Value = 16
$ws.Cells.Item(42, 7).Value = 16
$ws.Cells.Item(42, 8).Value = 15

$ws.Cells.Item(43, 1).Value = "Daniele❤"
$ws.Cells.Item(43, 2).Value = "Ok"
$ws.Cells.Item(43, 3).Value = "Guerra Atual"
$ws.Cells.Item(43, 4).Value = 16
$ws.Cells.Item(43, 5).Value = 16
$ws.Cells.Item(43, 6).Value = 16
$ws.Cells.Item(43, 7).Value = 16
$ws.Cells.Item(43, 8).Value = 16

$ws.Cells.Item(44, 1).Value = "GabiMalvadeza"
$ws.Cells.Item(44, 2).Value = "Razoável"
$ws.Cells.Item(44, 3).Value = "Guerra Atual"
$ws.Cells.Item(44, 4).Value = 15
$ws.Cells.Item(44, 5).Value = 0
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 0
$ws.Cells.Item(44, 8).Value = 0

$ws.Cells.Item(45, 1).Value = "StelaAby"
$ws.Cells.Item(45, 2).Value = "Verificar"
$ws.Cells.Item(45, 3).Value = "Guerra Atual"
$ws.Cells.Item(45, 4).Value = 8
$ws.Cells.Item(45, 5).Value = 0
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(45, 7).Value = 0
$ws.Cells.Item(45, 8).Value = 0

$ws.Cells.Item(46, 1).Value = "⭐O SENTINELA ⭐"
$ws.Cells.Item(46, 2).Value = "Ok"
$ws.Cells.Item(46, 3).Value = "Guerra Atual"
$ws.Cells.Item(46, 4).Value = 16
$ws.Cells.Item(46, 5).Value = 16
$ws.Cells.Item(46, 6).Value = 16
$ws.Cells.Item(46, 7).Value = 16
$ws.Cells.Item(46, 8).Value = 16

$ws.Cells.Item(47, 1).Value = "WvCly"
$ws.Cells.Item(47, 2).Value = "Ok"
$ws.Cells.Item(47, 3).Value = "Guerra Atual"
$ws.Cells.Item(47, 4).Value = 16
$ws.Cells.Item(47, 5).Value = 8
$ws.Cells.Item(47, 6).Value = 14
$ws.Cells.Item(47, 7).Value = 7
$ws.Cells.Item(47, 8).Value = 10

$ws.Cells.Item(48, 1).Value = "Pedro PH"
$ws.Cells.Item(48, 2).Value = "Razoável"
$ws.Cells.Item(48, 3).Value = "Guerra Atual"
$ws.Cells.Item(48, 4).Value = 15
$ws.Cells.Item(48, 5).Value = 16
$ws.Cells.Item(48, 6).Value = 16
$ws.Cells.Item(48, 7).Value = 8
$ws.Cells.Item(48, 8).Value = 0

$ws.Cells.Item(49, 1).Value = "Diih"
$ws.Cells.Item(49, 2).Value = "Verificar"
$ws.Cells.Item(49, 3).Value = "Guerra Atual"
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(49, 5).Value = 0
$ws.Cells.Item(49, 6).Value = 0
$ws.Cells.Item(49, 7).Value = 0
$ws.Cells.Item(49, 8).Value = 0

$ws.Cells.Item(50, 1).Value = "Chetto"
$ws.Cells.Item(50, 2).Value = "Razoável"
$ws.Cells.Item(50, 3).Value = "Guerra Atual"
$ws.Cells.Item(50, 4).Value = 14
$ws.Cells.Item(50, 5).Value = 12
$ws.Cells.Item(50, 6).Value = 13
$ws.Cells.Item(50, 7).Value = 11
$ws.Cells.Item(50, 8).Value = 16

$ws.Cells.Item(51, 1).Value = "Alvaro"
$ws.Cells.Item(51, 2).Value = "Razoável"
$ws.Cells.Item(51, 3).Value = "Guerra Atual"
$ws.Cells.Item(51, 4).Value = 15
$ws.Cells.Item(51, 5).Value = 13
$ws.Cells.Item(51, 6).Value = 14
$ws.Cells.Item(51, 7).Value = 15
$ws.Cells.Item(51, 8).Value = 14

$ws.Cells.Item(52, 1).Value = "domador de but"
$ws.Cells.Item(52, 2).Value = "Ok"
$ws.Cells.Item(52, 3).Value = "Guerra Atual"
$ws.Cells.Item(52, 4).Value = 16
$ws.Cells.Item(52, 5).Value = 16
$ws.Cells.Item(52, 6).Value = 16
$ws.Cells.Item(52, 7).Value = 8
$ws.Cells.Item(52, 8).Value = 12

$ws.Cells.Item(53, 1).Value = "Grimmer 狼"
$ws.Cells.Item(53, 2).Value = "Ok"
$ws.Cells.Item(53, 3).Value = "Guerra Atual"
$ws.Cells.Item(53, 4).Value = 16
$ws.Cells.Item(53, 5).Value = 16
$ws.Cells.Item(53, 6).Value = 4
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 0

$ws.Cells.Item(54, 1).Value = "gnvx v"
$ws.Cells.Item(54, 2).Value = "Verificar"
$ws.Cells.Item(54, 3).Value = "Guerra Atual"
$ws.Cells.Item(54, 4).Value = 0
$ws.Cells.Item(54, 5).Value = 0
$ws.Cells.Item(54, 6).Value = 0
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = 0

$ws.Cells.Item(55, 1).Value = "Nico"
$ws.Cells.Item(55, 2).Value = "Verificar"
$ws.Cells.Item(55, 3).Value = "Guerra Atual"
$ws.Cells.Item(55, 4).Value = 10
$ws.Cells.Item(55, 5).Value = 14
$ws.Cells.Item(55, 6).Value = 14
$ws.Cells.Item(55, 7).Value = 9
$ws.Cells.Item(55, 8).Value = 14

$ws.Cells.Item(56, 1).Value = "juniorllou"
$ws.Cells.Item(56, 2).Value = "Verificar"
$ws.Cells.Item(56, 3).Value = "Guerra Atual"
$ws.Cells.Item(56, 4).Value = 0
$ws.Cells.Item(56, 5).Value = 0
$ws.Cells.Item(56, 6).Value = 0
$ws.Cells.Item(56, 7).Value = 0
$ws.Cells.Item(56, 8).Value = 0

$ws.Cells.Item(57, 1).Value = "Fafazin10"
$ws.Cells.Item(57, 2).Value = "Verificar"
$ws.Cells.Item(57, 3).Value = "Guerra Atual"
$ws.Cells.Item(57, 4).Value = 0
$ws.Cells.Item(57, 5).Value = 0
$ws.Cells.Item(57, 6).Value = 0
$ws.Cells.Item(57, 7).Value = 0
$ws.Cells.Item(57, 8).Value = 0

$ws.Cells.Item(58, 1).Value = "Bruno"
$ws.Cells.Item(58, 2).Value = "Ok"
$ws.Cells.Item(58, 3).Value = "Guerra Atual"
$ws.Cells.Item(58, 4).Value = 16
$ws.Cells.Item(58, 5).Value = 16
$ws.Cells.Item(58, 6).Value = 16
$ws.Cells.Item(58, 7).Value = 0
$ws.Cells.Item(58, 8).Value = 0

$ws.Cells.Item(59, 1).Value = "Dockz"
$ws.Cells.Item(59, 2).Value = "Razoável"
$ws.Cells.Item(59, 3).Value = "Guerra Atual"
$ws.Cells.Item(59, 4).Value = 15
$ws.Cells.Item(59, 5).Value = 16
$ws.Cells.Item(59, 6).Value = 14
$ws.Cells.Item(59, 7).Value = 16
$ws.Cells.Item(59, 8).Value = 16

$ws.Cells.Item(60, 1).Value = "O GUARDIÃO"
$ws.Cells.Item(60, 2).Value = "Ok"
$ws.Cells.Item(60, 3).Value = "Guerra Atual"
$ws.Cells.Item(60, 4).Value = 16
$ws.Cells.Item(60, 5).Value = 16
$ws.Cells.Item(60, 6).Value = 16
$ws.Cells.Item(60, 7).Value = 16
$ws.Cells.Item(60, 8).Value = 16

$ws.Cells.Item(61, 1).Value = "Kauan"
$ws.Cells.Item(61, 2).Value = "Verificar"
$ws.Cells.Item(61, 3).Value = "Guerra Atual"
$ws.Cells.Item(61, 4).Value = 0
$ws.Cells.Item(61, 5).Value = 0
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(61, 8).Value = 0

$ws.Cells.Item(62, 1).Value = "Luciano"
$ws.Cells.Item(62, 2).Value = "Ok"
$ws.Cells.Item(62, 3).Value = "Guerra Atual"
$ws.Cells.Item(62, 4).Value = 16
$ws.Cells.Item(62, 5).Value = 16
$ws.Cells.Item(62, 6).Value = 16
$ws.Cells.Item(62, 7).Value = 14
$ws.Cells.Item(62, 8).Value = 16

$ws.Cells.Item(63, 1).Value = "Teixeirazzqw"
$ws.Cells.Item(63, 2).Value = "Ok"
$ws.Cells.Item(63, 3).Value = "Guerra Atual"
$ws.Cells.Item(63, 4).Value = 16
$ws.Cells.Item(63, 5).Value = 10
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 0

$ws.Cells.Item(64, 1).Value = "EDDIE"
$ws.Cells.Item(64, 2).Value = "Ok"
$ws.Cells.Item(64, 3).Value = "Guerra Atual"
$ws.Cells.Item(64, 4).Value = 16
$ws.Cells.Item(64, 5).Value = 16
$ws.Cells.Item(64, 6).Value = 16
$ws.Cells.Item(64, 7).Value = 16
$ws.Cells.Item(64, 8).Value = 16

$ws.Cells.Item(65, 1).Value = "isp"
$ws.Cells.Item(65, 2).Value = "Ok"
$ws.Cells.Item(65, 3).Value = "Guerra Atual"
$ws.Cells.Item(65, 4).Value = 16
$ws.Cells.Item(65, 5).Value = 16
$ws.Cells.Item(65, 6).Value = 16
$ws.Cells.Item(65, 7).Value = 16
$ws.Cells.Item(65, 8).Value = 16

$ws.Cells.Item(66, 1).Value = "filho de duque"
$ws.Cells.Item(66, 2).Value = "Ok"
$ws.Cells.Item(66, 3).Value = "Guerra Atual"
$ws.Cells.Item(66, 4).Value = 16
$ws.Cells.Item(66, 5).Value = 0
$ws.Cells.Item(66, 6).Value = 0
$ws.Cells.Item(66, 7).Value = 0
$ws.Cells.Item(66, 8).Value = 0

$ws.Cells.Item(67, 1).Value = "DGJ-DAVI"
$ws.Cells.Item(67, 2).Value = "Ok"
$ws.Cells.Item(67, 3).Value = "Guerra Atual"
$ws.Cells.Item(67, 4).Value = 16
$ws.Cells.Item(67, 5).Value = 16
$ws.Cells.Item(67, 6).Value = 16
$ws.Cells.Item(67, 7).Value = 12
$ws.Cells.Item(67, 8).Value = 15

$ws.Cells.Item(68, 1).Value = "luck"
$ws.Cells.Item(68, 2).Value = "Ok"
$ws.Cells.Item(68, 3).Value = "Guerra Atual"
$ws.Cells.Item(68, 4).Value = 16
$ws.Cells.Item(68, 5).Value = 16
$ws.Cells.Item(68, 6).Value = 16
$ws.Cells.Item(68, 7).Value = 12
$ws.Cells.Item(68, 8).Value = 16

$ws.Cells.Item(69, 1).Value = "andrebts"
$ws.Cells.Item(69, 2).Value = "Ok"
$ws.Cells.Item(69, 3).Value = "Guerra Atual"
$ws.Cells.Item(69, 4).Value = 16
$ws.Cells.Item(69, 5).Value = 15
$ws.Cells.Item(69, 6).Value = 16
$ws.Cells.Item(69, 7).Value = 16
$ws.Cells.Item(69, 8).Value = 16

$ws.Cells.Item(70, 1).Value = "Asten Acady"
$ws.Cells.Item(70, 2).Value = "Ok"
$ws.Cells.Item(70, 3).Value = "Guerra Atual"
$ws.Cells.Item(70, 4).Value = 16
$ws.Cells.Item(70, 5).Value = 16
$ws.Cells.Item(70, 6).Value = 16
$ws.Cells.Item(70, 7).Value = 16
$ws.Cells.Item(70, 8).Value = 16

$ws.Cells.Item(71, 1).Value = "OneDePrata"
$ws.Cells.Item(71, 2).Value = "Ok"
$ws.Cells.Item(71, 3).Value = "Guerra Atual"
$ws.Cells.Item(71, 4).Value = 16
$ws.Cells.Item(71, 5).Value = 15
$ws.Cells.Item(71, 6).Value = 8
$ws.Cells.Item(71, 7).Value = 16
$ws.Cells.Item(71, 8).Value = 13

$ws.Cells.Item(72, 1).Value = "51 é pinga"
$ws.Cells.Item(72, 2).Value = "Ok"
$ws.Cells.Item(72, 3).Value = "Guerra Atual"
$ws.Cells.Item(72, 4).Value = 16
$ws.Cells.Item(72, 5).Value = 16
$ws.Cells.Item(72, 6).Value = 16
$ws.Cells.Item(72, 7).Value = 16
$ws.Cells.Item(72, 8).Value = 16

$ws.Cells.Item(73, 1).Value = "gabiggoughost"
$ws.Cells.Item(73, 2).Value = "Verificar"
$ws.Cells.Item(73, 3).Value = "Guerra Atual"
$ws.Cells.Item(73, 4).Value = 8
$ws.Cells.Item(73, 5).Value = 16
$ws.Cells.Item(73, 6).Value = 14
$ws.Cells.Item(73, 7).Value = 16
$ws.Cells.Item(73, 8).Value = 16

$ws.Cells.Item(74, 1).Value = "nivelador"
$ws.Cells.Item(74, 2).Value = "Ok"
$ws.Cells.Item(74, 3).Value = "Guerra Atual"
$ws.Cells.Item(74, 4).Value = 16
$ws.Cells.Item(74, 5).Value = 12
$ws.Cells.Item(74, 6).Value = 10
$ws.Cells.Item(74, 7).Value = 15
$ws.Cells.Item(74, 8).Value = 12

$ws.Cells.Item(75, 1).Value = "joão3:16"
$ws.Cells.Item(75, 2).Value = "Ok"
$ws.Cells.Item(75, 3).Value = "Guerra Atual"
$ws.Cells.Item(75, 4).Value = 16
$ws.Cells.Item(75, 5).Value = 16
$ws.Cells.Item(75, 6).Value = 16
$ws.Cells.Item(75, 7).Value = 16
$ws.Cells.Item(75, 8).Value = 16

$ws.Cells.Item(76, 1).Value = "Mila"
$ws.Cells.Item(76, 2).Value = "Ok"
$ws.Cells.Item(76, 3).Value = "Guerra Atual"
$ws.Cells.Item(76, 4).Value = 16
$ws.Cells.Item(76, 5).Value = 8
$ws.Cells.Item(76, 6).Value = 16
$ws.Cells.Item(76, 7).Value = 16
$ws.Cells.Item(76, 8).Value = 16

$ws.Cells.Item(77, 1).Value = "Sotto ツ"
$ws.Cells.Item(77, 2).Value = "Razoável"
$ws.Cells.Item(77, 3).Value = "Guerra Atual"
$ws.Cells.Item(77, 4).Value = 12
$ws.Cells.Item(77, 5).Value = 12
$ws.Cells.Item(77, 6).Value = 14
$ws.Cells.Item(77, 7).Value = 12
$ws.Cells.Item(77, 8).Value = 16
